$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Unidad de Asignacion" / "Notas" columns (D, E) entirely
$ws.Range("D1:E7").Delete()

# Clear the old data rows (keep header row intact, it gets edited in place below)
$ws.Range("A2:C7").Clear()

# Update header row text (B1, C1 renamed; A1 stays the same)
$ws.Range("B1").Value = "Tamaño del Archivo"
$ws.Range("C1").Value = "Tiempo de Transferencia"

# Fill column A (Sistema de Archivos)
$ws.Range("A2").Value = "exFAT"
$ws.Range("A3").Value = "exFAT"
$ws.Range("A4").Value = "exFAT"
$ws.Range("A5").Value = "NTFS"
$ws.Range("A6").Value = "NTFS"
$ws.Range("A7").Value = "NTFS"

# Fill column C (Tiempo de Transferencia)
$ws.Range("C2").Value = "1 seg"
$ws.Range("C3").Value = "10 seg"
$ws.Range("C4").Value = "11 min"
$ws.Range("C5").Value = "1 seg"
$ws.Range("C6").Value = "10 seg"
$ws.Range("C7").Value = "2 min 15 seg"

# Fill column B (Tamaño del Archivo)
$ws.Range("B2").Value = "500 KB"
$ws.Range("B3").Value = "500 MB"
$ws.Range("B4").Value = "1 GB"
$ws.Range("B5").Value = "500 KB"
$ws.Range("B6").Value = "500 MB"
$ws.Range("B7").Value = "1 GB"

# Data rows style: thin border only (rows 2-7)
$dataRange = $ws.Range("A2:C7")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# Header style: bold font, centered horizontal/top vertical alignment, thin border
$headerRange = $ws.Range("A1:C1")
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Column widths to match the new, narrower 3-column layout
$ws.Columns.Item(1).ColumnWidth = 23.5546875
$ws.Columns.Item(2).ColumnWidth = 25.33203125
$ws.Columns.Item(3).ColumnWidth = 22.44140625

# Update the selection to match the target workbook state
$ws.Range("F9").Select()
